$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 85 (was the 43979 row, becomes the 43978 row)
$ws.Range("B85").Value = 43978.0
$ws.Range("C85").Value = 2140.0
$ws.Range("D85").Value = 171449.0
$ws.Range("E85").Value = 4.0

# Fix row 86 (was the 43978 row, becomes the 43979 row)
$ws.Range("B86").Value = 43979.0
$ws.Range("C86").Value = 2054.0
$ws.Range("D86").Value = 173503.0
$ws.Range("E86").Value = 5.0

# Add new row 87
# Force A87 to be stored as text (matching the "weekday counter as text"
# pattern used by the rest of column A), then strip the formatting change
# back off so no stray number format lingers on the cell.
$ws.Range("A87").NumberFormat = "@"
$ws.Range("A87").Value = "86"
$ws.Range("A87").ClearFormats()

$ws.Range("B87").Value = 43980.0
$ws.Range("C87").Value = 2332.0
$ws.Range("D87").Value = 175835.0
$ws.Range("E87").Value = 6.0

# Make sure B87 has the same date number format as the other date cells
$ws.Range("B87").NumberFormat = $ws.Range("B86").NumberFormat
